$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# Update workbook window height (book view)
$excel.ActiveWindow.Height = 16440

# Update the sheet's selected range (RP all / new McT algorithm block)
$ws.Range("J1:L1").Select()

# J3 / K3
$ws.Range("J3").Value = 1.1100000000000001
$ws.Range("K3").Formula = "=0.44/2"

# J4 / K4
$ws.Range("J4").Value = 1.01
$ws.Range("K4").Formula = "=0.32/2"

# J5 / K5
$ws.Range("J5").Value = 1.25
$ws.Range("K5").Formula = "=0.53/2"

# K6
$ws.Range("K6").Formula = "=0.42/2"

# J8 / K8
$ws.Range("J8").Value = 1.03
$ws.Range("K8").Formula = "=0.63/2"

# K9
$ws.Range("K9").Formula = "=0.56/2"

# K10 / K11 -> clear values
$ws.Range("K10").ClearContents()
$ws.Range("K11").ClearContents()

# K13
$ws.Range("K13").Formula = "=1.29/2"

# J15 / K15
$ws.Range("J15").Value = 0.99
$ws.Range("K15").Formula = "=0.99/2"
